# Generate Report for Handoff
#
# The localization status report previously showed "b.md" as already handed
# back (in sync with en-US). A newer handback package for b.md has now been
# processed, so its status flips to "Ready for handoff" (it is not yet fully
# in sync - the handback file version does not match the latest source), and
# the related handoff/handback metadata + error detail are refreshed across
# all three report sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (b.md) status + date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 08:36:43"

# --- zh-cn sheet: row 3 (b.md) status/content-duplicate/handback-file/date/error ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces these as literal text (matching the source file,
# where "True"/"False" are stored as plain strings, not Excel booleans).
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-17 08:36:38"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66f9e7aa1dd82195b8c52fecac6edc390cac4711/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a8d9273838290472fb9c081318ca8f40f961b8e/e2e/b.md."

# --- de-de sheet: row 3 (b.md) status/content-duplicate/handback-file/date/error ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-17 08:36:43"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66f9e7aa1dd82195b8c52fecac6edc390cac4711/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a8d9273838290472fb9c081318ca8f40f961b8e/e2e/b.md."

# --- Column P (Error Detail) is widened to fit the new, longer message text ---
$wsZhCn.Columns("P").ColumnWidth = 40
$wsDeDe.Columns("P").ColumnWidth = 40
